# Apply updated "dSF" (column F) values to Sheet1, rows 2-30.
# These reflect a data repull / recalculation of the dSF metric.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 1
    5  = 12
    7  = 1
    8  = -1
    9  = -1
    10 = 3
    12 = 3
    13 = 4
    14 = 2
    15 = 3
    16 = 3
    17 = -3
    18 = 5
    19 = 5
    20 = 1
    21 = 1
    22 = 1
    24 = 5
    25 = -1
    26 = 4
    27 = -3
    28 = -1
    29 = 6
    30 = -1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 6).Value = $newValues[$row]
}
